$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valid_Login")

# Update the password cell C6 with the new value and turn it into a mailto-style
# hyperlink, matching the pattern used by the other password cell (C5).
$ws.Range("C6").Value = "Surekha@123"
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:Surekha@123") | Out-Null
$ws.Range("C6").Style = "Hyperlink"
